# RCC new script implementation
# Adds two new test-case rows (RCC115, RCC116) to the "Test Cases" sheet,
# mirroring the formatting of the existing row 18 (RCC114).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 19 (RCC115): clone formatting from row 18 ---------------------
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)

# --- Row 20 (RCC116): clone formatting from row 18 ----------------------
$ws.Range("A18:E18").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)

# Column C of row 18 (style index 6) isn't what row 19 needs for its
# Description cell - row 19's Description cell matches the "wrap + border"
# style used by column B of row 18, so copy that style across instead.
$ws.Range("B18").Copy()
$ws.Range("C19").PasteSpecial(-4122)

# Row 20's Jira-id cell (column B) drops its border and keeps only wrap
# text - remove the border and ensure wrap text is applied.
$ws.Range("B20").Borders.LineStyle = 0
$ws.Range("B20").WrapText = $true

# Both new rows are tall (wrapped, multi-line content) like row 18.
$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 45

# --- Values for row 19 (RCC115) -----------------------------------------
$ws.Range("A19").Value = "RCC115"
$ws.Range("C19").Value = "Verify that user is able to add an article to the multiple groups from record view page.||Verify that user is able to add a post to the  multiple groups from record view page. ||Verify that user is able to add a patent to the  multiple groups from record view page."
$ws.Range("B19").Value = "OPQA-3467||OPQA-3471||OPQA-3475"
$ws.Range("D19").Value = "Y"

# --- Values for row 20 (RCC116) -----------------------------------------
$ws.Range("A20").Value = "RCC116"
$ws.Range("C20").Value = "Verify that user is able to add an article to the multiple groups from watch list details page.||Verify that user is able to add a post to the  multiple groups from watch list details page.||Verify that user is able to add a patent to the  multiple groups from watch list details page."
$ws.Range("B20").Value = "OPQA-3469||OPQA-3473||OPQA-3477"
$ws.Range("D20").Value = "Y"

# Match the final selection/active cell left by the author's edit.
$ws.Range("B20").Select()
